$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: Lauri Markkanen, SF,PF, Utah Jazz  ->  RJ Barrett, SF,PF, Toronto Raptors
$ws.Range("A6").Value = "RJ Barrett"
$ws.Range("C6").Value = "Toronto Raptors"

# Row 10: John Collins, PF,C, Utah Jazz  ->  Dorian Finney-Smith, SF,PF,C, Brooklyn Nets
$ws.Range("A10").Value = "Dorian Finney-Smith"
$ws.Range("B10").Value = "SF,PF,C"
$ws.Range("C10").Value = "Brooklyn Nets"

# Row 11: Dorian Finney-Smith, SF,PF,C, Brooklyn Nets  ->  Joel Embiid, C, Philadelphia 76ers
$ws.Range("A11").Value = "Joel Embiid"
$ws.Range("B11").Value = "C"
$ws.Range("C11").Value = "Philadelphia 76ers"

# Row 12: Jalen Williams, SG,SF,PF, Oklahoma City Thunder  ->  Jalen Williams, SG,SF,PF,C, Oklahoma City Thunder
$ws.Range("B12").Value = "SG,SF,PF,C"

# Row 14: Joel Embiid, C, Philadelphia 76ers  ->  Lauri Markkanen, SF,PF, Utah Jazz
$ws.Range("A14").Value = "Lauri Markkanen"
$ws.Range("B14").Value = "SF,PF"
$ws.Range("C14").Value = "Utah Jazz"

# Row 15: RJ Barrett, SF,PF, Toronto Raptors  ->  John Collins, PF,C, Utah Jazz
$ws.Range("A15").Value = "John Collins"
$ws.Range("B15").Value = "PF,C"
$ws.Range("C15").Value = "Utah Jazz"
